{"js": "// The edit rotates the text of eight runs through the document:\n// Objetivos -> Docente -> Bibliografia, and\n// Programa resumido -> Programa -> M\u00e9todo -> Crit\u00e9rio -> Norma de recupera\u00e7\u00e3o -> Programa resumido.\n// Concretely (old text -> new text):\n//   Objetivos:            \"Falar basicamente...\"      -> \"Ci\u00eancia ontem e hoje...\"\n//   Docente:               \"6376612 - Daisy...\"        -> \"Falar basicamente...\"\n//   Programa resumido:     \"Ci\u00eancia ontem e hoje...\"    -> \"1. Ci\u00eancia, t\u00e9cnica...\"\n//   Programa:              \"1. Ci\u00eancia, t\u00e9cnica...\"     -> \"Aulas expositivas...\"\n//   M\u00e9todo:                \"Aulas expositivas...\"       -> \"A = (P + T)/ 2...\"\n//   Crit\u00e9rio:               \"A = (P + T)/ 2...\"          -> \"RECUPERA\u00c7\u00c3O: 1 (uma) prova.\"\n//   Norma de recupera\u00e7\u00e3o:  \"RECUPERA\u00c7\u00c3O: 1 (uma) prova.\" -> \"VARGAS, Milton...\"\n//   Bibliografia:           \"VARGAS, Milton...\"          -> \"6376612 - Daisy...\"\n//\n// We look each old string up with body.search (unique matches), capture the\n// matching ranges first, then assign new text to each -- this avoids any\n// ordering hazard from one replacement's new text accidentally matching\n// another block's search term.\n\nconst oldTexts = [\n  \"Falar basicamente sobre ci\u00eancia; T\u00e9cnica; Tecnologia; Engenharia; Pesquisa; Descobertas e Inven\u00e7\u00f5es.\",\n  \"6376612 - Daisy Rafaela da Silva\",\n  \"Ci\u00eancia ontem e hoje; Positivismo e o Neopositivismo; F\u00edsica moderna e seus pensadores; As escolas de Engenharia; Prepara\u00e7\u00e3o de monografias.\",\n  \"1. Ci\u00eancia, t\u00e9cnica, tecnologia e engenharia 2. Ci\u00eancia e t\u00e9cnica na Idade Antiga 3. Ci\u00eancia e t\u00e9cnica na Idade M\u00e9dia . 4. Ci\u00eancia e t\u00e9cnica na Idade Moderna. 5. Ci\u00eancia e t\u00e9cnica na Idade Contempor\u00e2nea 6. Metodologia Cient\u00edfica: Arist\u00f3teles, Galileu e Descartes. 7. Positivismo e neopositivismo, COMTE, Popper, Kuhn. 8. Defini\u00e7\u00e3o, medidas, leis e teoria f\u00edsica 9. F\u00edsica Moderna e realidade. 10. Teoria dos Modelos. Newton, Poincar\u00e9, Lorentz, Einstein 11. Pesquisa de causas. Leis estat\u00edsticas, determinismo e acaso 12. Ci\u00eancia te\u00f3rica e ci\u00eancia experimental 13. Pesquisas, descobertas e inven\u00e7\u00f5es 14. As escolas de engenharia, forma\u00e7\u00e3o das escolas, Escola de Engenharia de Lorena /EEL/USP. 15. Engenharia, matem\u00e1tica e f\u00edsica 16. Organiza\u00e7\u00e3o da pesquisa tecnol\u00f3gica 17. Prepara\u00e7\u00e3o de Monografias tecnol\u00f3gicas.\",\n  \"Aulas expositivas em n\u00edvel de confer\u00eancia. Estudo de casos significativos da hist\u00f3ria da ci\u00eancia e da engenharia. Debate participativo em torno de quest\u00f5es relevantes.\",\n  \"A = (P + T)/ 2 Onde: P = m\u00e9dia das provas T = m\u00e9dia dos trabalhos pr\u00e1ticos\",\n  \"RECUPERA\u00c7\u00c3O: 1 (uma) prova.\",\n  \"VARGAS, Milton. Metodologia da pesquisa tecnol\u00f3gica,Rio de Janeiro, Globo, 1985. SIMARD, Emile. Naturaleza y alcance del m\u00e9todo cient\u00edfico, Madrim, Gredos, 1961. ROUSSEAU, Pierre. Hist\u00f3ria da Ci\u00eancia, Lisboa, 1963. VARGAS, Milton. Para uma filosofia da tecnologia, Ed.Alfa-Omega, S\u00e3o Paulo, 1994.\"\n];\n\n// new text for slot i is the old text that used to occupy the *next* slot\n// (wrapping around), matching the diff exactly.\nconst newTexts = [\n  oldTexts[2], // Objetivos            <- old Programa resumido\n  oldTexts[0], // Docente              <- old Objetivos\n  oldTexts[3], // Programa resumido    <- old Programa\n  oldTexts[4], // Programa             <- old M\u00e9todo\n  oldTexts[5], // M\u00e9todo               <- old Crit\u00e9rio\n  oldTexts[6], // Crit\u00e9rio             <- old Norma de recupera\u00e7\u00e3o\n  oldTexts[7], // Norma de recupera\u00e7\u00e3o <- old Bibliografia\n  oldTexts[1]  // Bibliografia         <- old Docente\n];\n\nconst body = context.document.body;\nconst ranges = [];\nfor (const t of oldTexts) {\n  const found = body.search(t, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n  ranges.push(found.items[0]);\n}\n\nfor (let i = 0; i < ranges.length; i++) {\n  ranges[i].insertText(newTexts[i], \"Replace\");\n}\nawait context.sync();\n", "ps1": "# The edit rotates the text of eight runs through the document:\n#   Objetivos:             \"Falar basicamente...\"       -> \"Ci\u00eancia ontem e hoje...\"\n#   Docente:                \"6376612 - Daisy...\"         -> \"Falar basicamente...\"\n#   Programa resumido:      \"Ci\u00eancia ontem e hoje...\"     -> \"1. Ci\u00eancia, t\u00e9cnica...\"\n#   Programa:               \"1. Ci\u00eancia, t\u00e9cnica...\"      -> \"Aulas expositivas...\"\n#   M\u00e9todo:                 \"Aulas expositivas...\"        -> \"A = (P + T)/ 2...\"\n#   Crit\u00e9rio:                \"A = (P + T)/ 2...\"           -> \"RECUPERA\u00c7\u00c3O: 1 (uma) prova.\"\n#   Norma de recupera\u00e7\u00e3o:   \"RECUPERA\u00c7\u00c3O: 1 (uma) prova.\" -> \"VARGAS, Milton...\"\n#   Bibliografia:            \"VARGAS, Milton...\"           -> \"6376612 - Daisy...\"\n#\n# Every \"new\" string here is some other block's \"old\" string, so a naive\n# sequential find/replace risks a later search accidentally matching text\n# that an earlier step just wrote. To avoid that, we first locate every\n# match (by character offsets, while the document still holds the original\n# text throughout), then apply the writes from the end of the document\n# backwards -- editing a later span can never shift the offsets of spans\n# that start earlier, so every captured offset stays valid until it's used.\n\n$d = $word.ActiveDocument\n\nfunction Find-ExactSpan($doc, $text) {\n    $rng = $doc.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $text\n    $null = $rng.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false)\n    return @($rng.Start, $rng.End)\n}\n\n$oldTexts = @(\n    \"Falar basicamente sobre ci\u00eancia; T\u00e9cnica; Tecnologia; Engenharia; Pesquisa; Descobertas e Inven\u00e7\u00f5es.\",\n    \"6376612 - Daisy Rafaela da Silva\",\n    \"Ci\u00eancia ontem e hoje; Positivismo e o Neopositivismo; F\u00edsica moderna e seus pensadores; As escolas de Engenharia; Prepara\u00e7\u00e3o de monografias.\",\n    \"1. Ci\u00eancia, t\u00e9cnica, tecnologia e engenharia 2. Ci\u00eancia e t\u00e9cnica na Idade Antiga 3. Ci\u00eancia e t\u00e9cnica na Idade M\u00e9dia . 4. Ci\u00eancia e t\u00e9cnica na Idade Moderna. 5. Ci\u00eancia e t\u00e9cnica na Idade Contempor\u00e2nea 6. Metodologia Cient\u00edfica: Arist\u00f3teles, Galileu e Descartes. 7. Positivismo e neopositivismo, COMTE, Popper, Kuhn. 8. Defini\u00e7\u00e3o, medidas, leis e teoria f\u00edsica 9. F\u00edsica Moderna e realidade. 10. Teoria dos Modelos. Newton, Poincar\u00e9, Lorentz, Einstein 11. Pesquisa de causas. Leis estat\u00edsticas, determinismo e acaso 12. Ci\u00eancia te\u00f3rica e ci\u00eancia experimental 13. Pesquisas, descobertas e inven\u00e7\u00f5es 14. As escolas de engenharia, forma\u00e7\u00e3o das escolas, Escola de Engenharia de Lorena /EEL/USP. 15. Engenharia, matem\u00e1tica e f\u00edsica 16. Organiza\u00e7\u00e3o da pesquisa tecnol\u00f3gica 17. Prepara\u00e7\u00e3o de Monografias tecnol\u00f3gicas.\",\n    \"Aulas expositivas em n\u00edvel de confer\u00eancia. Estudo de casos significativos da hist\u00f3ria da ci\u00eancia e da engenharia. Debate participativo em torno de quest\u00f5es relevantes.\",\n    \"A = (P + T)/ 2 Onde: P = m\u00e9dia das provas T = m\u00e9dia dos trabalhos pr\u00e1ticos\",\n    \"RECUPERA\u00c7\u00c3O: 1 (uma) prova.\",\n    \"VARGAS, Milton. Metodologia da pesquisa tecnol\u00f3gica,Rio de Janeiro, Globo, 1985. SIMARD, Emile. Naturaleza y alcance del m\u00e9todo cient\u00edfico, Madrim, Gredos, 1961. ROUSSEAU, Pierre. Hist\u00f3ria da Ci\u00eancia, Lisboa, 1963. VARGAS, Milton. Para uma filosofia da tecnologia, Ed.Alfa-Omega, S\u00e3o Paulo, 1994.\"\n)\n\n# new text for slot i is the old text that used to occupy the *next* slot\n# (wrapping around), matching the diff exactly.\n$newTexts = @(\n    $oldTexts[2], # Objetivos            <- old Programa resumido\n    $oldTexts[0], # Docente              <- old Objetivos\n    $oldTexts[3], # Programa resumido    <- old Programa\n    $oldTexts[4], # Programa             <- old M\u00e9todo\n    $oldTexts[5], # M\u00e9todo               <- old Crit\u00e9rio\n    $oldTexts[6], # Crit\u00e9rio             <- old Norma de recupera\u00e7\u00e3o\n    $oldTexts[7], # Norma de recupera\u00e7\u00e3o <- old Bibliografia\n    $oldTexts[1]  # Bibliografia         <- old Docente\n)\n\n$spans = @()\nfor ($i = 0; $i -lt $oldTexts.Length; $i++) {\n    $span = Find-ExactSpan $d $oldTexts[$i]\n    $spans += , @($i, $span[0], $span[1])\n}\n\n$order = $spans | Sort-Object { $_[1] } -Descending\n\nforeach ($entry in $order) {\n    $idx = $entry[0]\n    $rng = $d.Range($entry[1], $entry[2])\n    $rng.Text = $newTexts[$idx]\n}\n"}
